$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 3-6 with new gene data (values rotate from old row 6 -> 3, 3 -> 4, 4 -> 5, 5 -> 6)
$ws.Cells.Item(3, 1).Value = "Rv0669c"
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = "Rv0669c"
$ws.Cells.Item(3, 4).Value = "FUNCTION: Catalyzes the cleavage of the N-acyl linkage of the ceramides (Cers) to yield sphingosine (Sph) and free fatty acid. Also catalyzes the synthesis of Cers from Sph and fatty acid. Cers containning C6-C24 fatty acids are well hydrolyzed, and Cers with mono unsaturated fatty acids are much more hydrolyzed than those with saturated fatty acids. {ECO:0000269|PubMed:10593963, ECO:0000269|PubMed:20139604}."
$ws.Cells.Item(3, 5).Value = 24

$ws.Cells.Item(4, 1).Value = "Rv0451c"
$ws.Cells.Item(4, 2).Value = 4
$ws.Cells.Item(4, 3).Value = "mmpS4 Rv0451c MTV037.15c"
$ws.Cells.Item(4, 4).Value = "FUNCTION: Part of an export system, which is required for biosynthesis and secretion of siderophores. Essential for virulence. {ECO:0000269|PubMed:23431276}."
$ws.Cells.Item(4, 5).Value = 24

$ws.Cells.Item(5, 1).Value = "Rv3537"
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = "kstD Rv3537"
$ws.Cells.Item(5, 4).Value = "FUNCTION: Involved in the degradation of cholesterol (PubMed:18031290, PubMed:21987574). Catalyzes the elimination of the C-1 and C-2 hydrogen atoms of the A-ring from the polycyclic ring structure of 3-ketosteroids (PubMed:18031290). Has a clear preference for 3-ketosteroids with a saturated A-ring, displaying highest activity on 5alpha-AD (5alpha-androstane-3,17-dione) and 5alpha-T (5alpha-testosterone, also known as 17beta-hydroxy-5alpha-androstane-3-one) (PubMed:18031290). Is also involved in the formation of 3-keto-1,4-diene-steroid from 3-keto-4-ene-steroid (PubMed:21987574). Catalyzes the conversion of 3-oxo-23,24-bisnorchol-4-en-22-oyl-coenzyme A thioester (4-BNC-CoA) to 3-oxo-23,24-bisnorchola-1,4-dien-22-oyl-coenzyme A thioester (1,4-BNC-CoA) (PubMed:21987574). {ECO:0000269|PubMed:18031290, ECO:0000269|PubMed:21987574}."
$ws.Cells.Item(5, 5).Value = 24

$ws.Cells.Item(6, 1).Value = "Rv2476c"
$ws.Cells.Item(6, 2).Value = 4
$ws.Cells.Item(6, 3).Value = "gdh Rv2476c"
$ws.Cells.Item(6, 4).Value = "FUNCTION: Catalyzes the reversible conversion of L-glutamate to 2-oxoglutarate. {ECO:0000250}."
$ws.Cells.Item(6, 5).Value = 24

# Append new rows 8-21 with additional gene-cluster data
$ws.Cells.Item(8, 1).Value = "Rv3553"
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = "ipdC Rv3553"
$ws.Cells.Item(8, 4).Value = "FUNCTION: Involved in the final steps of cholesterol and steroid degradation (PubMed:28377529). Probably catalyzes the introduction of a double bound into the C ring of 5OH-HIC-CoA, leading to the formation of (5R,7aS)-5-hydroxy-7a-methyl-1-oxo-3,5,6,7-tetrahydro-2H-indene-4-carboxyl-CoA (Probable). {ECO:0000269|PubMed:28377529, ECO:0000305|PubMed:28377529}."
$ws.Cells.Item(8, 5).Value = 24

$ws.Cells.Item(9, 1).Value = "Rv0693"
$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(9, 3).Value = "mftC Rv0693"
$ws.Cells.Item(9, 4).Value = "FUNCTION: Radical S-adenosylmethionine (SAM) enzyme responsible for the first step of the biosynthesis of the enzyme cofactor mycofactocin (MFT). Catalyzes two reactions at the C-terminus of the mycofactocin precursor (the MftA peptide). The first one is the oxidative decarboxylation of the C-terminal L-tyrosine of MftA, forming an unsaturated tyramine moiety. The second reaction is the cross-linking of the tyramine with the penultimate L-valine residue, forming a five-membered lactam ring. Its activity requires the presence of the MftB chaperone. {ECO:0000250|UniProtKB:A0PM49}."
$ws.Cells.Item(9, 5).Value = 24

$ws.Cells.Item(10, 1).Value = "Rv0391"
$ws.Cells.Item(10, 2).Value = 3
$ws.Cells.Item(10, 3).Value = "metZ Rv0391"
$ws.Cells.Item(10, 4).Value = "FUNCTION: Catalyzes the formation of L-homocysteine from O-succinyl-L-homoserine (OSHS) and hydrogen sulfide. {ECO:0000255|HAMAP-Rule:MF_02056}."
$ws.Cells.Item(10, 5).Value = 24

$ws.Cells.Item(11, 1).Value = "Rv1273c"
$ws.Cells.Item(11, 2).Value = 2
$ws.Cells.Item(11, 3).Value = "Rv1273c MTCY50.09"
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = 24

$ws.Cells.Item(12, 1).Value = "Rv0362"
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(12, 3).Value = "mgtE Rv0362"
$ws.Cells.Item(12, 4).Value = "FUNCTION: Acts as a magnesium transporter. {ECO:0000256|RuleBase:RU362011}."
$ws.Cells.Item(12, 5).Value = 24

$ws.Cells.Item(13, 1).Value = "Rv1127c"
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = "ppdK Rv1127c"
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = 24

$ws.Cells.Item(14, 1).Value = "Rv3493c"
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = "Rv3493c"
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = 24

$ws.Cells.Item(15, 1).Value = "Rv3548c"
$ws.Cells.Item(15, 2).Value = 1
$ws.Cells.Item(15, 3).Value = "Rv3548c"
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = 24

$ws.Cells.Item(16, 1).Value = "Rv3549c"
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(16, 3).Value = "Rv3549c"
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 5).Value = 24

$ws.Cells.Item(17, 1).Value = "Rv2047c"
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = "Rv2047c"
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = 24

$ws.Cells.Item(18, 1).Value = "Rv1627c"
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = "Rv1627c"
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(18, 5).Value = 24

$ws.Cells.Item(19, 1).Value = "Rv3503c"
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(19, 3).Value = "fdxD Rv3503c"
$ws.Cells.Item(19, 4).Value = "FUNCTION: Ferredoxins are iron-sulfur proteins that transfer electrons in a wide variety of metabolic reactions. {ECO:0000256|RuleBase:RU368020}."
$ws.Cells.Item(19, 5).Value = 24

$ws.Cells.Item(20, 1).Value = "Rv0390"
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = "Rv0390"
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = 24

$ws.Cells.Item(21, 1).Value = "Rv0320"
$ws.Cells.Item(21, 2).Value = 1
$ws.Cells.Item(21, 3).Value = "Rv0320"
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(21, 5).Value = 24
